$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A80").Value = "X"
$ws.Range("A76:D76").Copy()
$ws.Range("A80").PasteSpecial(-4122)
